$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header B1: "dihadron" -> "hadrons"
$ws.Range("B1").Value = "hadrons"

# Data column B rows 2-59: "pion" -> "2(pi+,pi-)"
$rng = $ws.Range("B2:B59")
$rng.Value = "2(pi+,pi-)"

# Update the active selection to L2
$ws.Range("L2").Select()
